# Refresh the dataset: multi-index (Date/Channel/Metric) rows feeding the upcoming
# multi-index/multi-table helpers in the gen-datasets script. Column A mixes real dates
# with the two sentinel text values ("yesterday" / "not a date") the new parsing code expects.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old body rows (A2:D26) so stale cells/styles from the previous 25-row dataset
# do not linger once the new 20-row dataset (A2:D21) is written.
$ws.Range("A2:D26").Clear()

# Row 2
$ws.Range("A2").Value = 45961
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B2").Value = "TV"
$ws.Range("C2").Value = "Spend"

# Row 3
$ws.Range("B3").Value = "TV"
$ws.Range("C3").Value = "GRPs"
$ws.Range("D3").Value = 5

# Row 4
$ws.Range("A4").Value = 45982
$ws.Range("A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B4").Value = "Radio"
$ws.Range("C4").Value = "GRPs"
$ws.Range("D4").Value = 2

# Row 5
$ws.Range("A5").Value = 45961
$ws.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B5").Value = "TV"
$ws.Range("C5").Value = "GRPs"
$ws.Range("D5").Value = 8

# Row 6
$ws.Range("A6").Value = 45975
$ws.Range("A6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B6").Value = "TV"
$ws.Range("C6").Value = "Spend"

# Row 7
$ws.Range("A7").Value = 45968
$ws.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B7").Value = "TV"
$ws.Range("C7").Value = "GRPs"
$ws.Range("D7").Value = 10

# Row 8
$ws.Range("A8").Value = "yesterday"
$ws.Range("B8").Value = "Radio"
$ws.Range("C8").Value = "GRPs"
$ws.Range("D8").Value = 10

# Row 9
$ws.Range("A9").Value = 45961
$ws.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B9").Value = "Radio"
$ws.Range("C9").Value = "GRPs"
$ws.Range("D9").Value = 6

# Row 10
$ws.Range("A10").Value = 45989
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B10").Value = "Radio"
$ws.Range("C10").Value = "Spend"
$ws.Range("D10").Value = 198

# Row 11
$ws.Range("A11").Value = 45989
$ws.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B11").Value = "TV"
$ws.Range("C11").Value = "Spend"
$ws.Range("D11").Value = 68

# Row 12
$ws.Range("A12").Value = 45982
$ws.Range("A12").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B12").Value = "TV"
$ws.Range("C12").Value = "GRPs"
$ws.Range("D12").Value = 4

# Row 13
$ws.Range("A13").Value = 45961
$ws.Range("A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B13").Value = "Radio"
$ws.Range("C13").Value = "Spend"
$ws.Range("D13").Value = "NaN"

# Row 14
$ws.Range("A14").Value = 45975
$ws.Range("A14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B14").Value = "TV"
$ws.Range("C14").Value = "GRPs"
$ws.Range("D14").Value = 8

# Row 15
$ws.Range("A15").Value = 45989
$ws.Range("A15").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B15").Value = "Radio"
$ws.Range("C15").Value = "GRPs"
$ws.Range("D15").Value = 4

# Row 16
$ws.Range("A16").Value = 45968
$ws.Range("A16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B16").Value = "TV"
$ws.Range("C16").Value = "Spend"
$ws.Range("D16").Value = 89

# Row 17
$ws.Range("A17").Value = 45982
$ws.Range("A17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B17").Value = "TV"
$ws.Range("C17").Value = "Spend"
$ws.Range("D17").Value = 110

# Row 18
$ws.Range("A18").Value = 45968
$ws.Range("A18").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B18").Value = "Radio"
$ws.Range("C18").Value = "GRPs"
$ws.Range("D18").Value = "NaN"

# Row 19
$ws.Range("A19").Value = 45975
$ws.Range("A19").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B19").Value = "Radio"
$ws.Range("C19").Value = "Spend"

# Row 20
$ws.Range("A20").Value = "not a date"
$ws.Range("B20").Value = "Radio"
$ws.Range("C20").Value = "Spend"
$ws.Range("D20").Value = 62

# Row 21
$ws.Range("A21").Value = 45968
$ws.Range("A21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B21").Value = "Radio"
$ws.Range("C21").Value = "Spend"
$ws.Range("D21").Value = 51
